# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-18 08:22:04
#
# Applies the attendance-report refresh to the "Session Analysis Results" sheet:
#   - swap the "Recorded By" ordering for the System/dnasr281 rows
#   - narrow the "Students" column
#   - refresh summary counters (Recorded/Missing sessions, coverage %, attendance %)
#   - refresh the per-group breakdown numbers (O/P/R/S columns)
#   - mark previously "Not Recorded" sessions (rows 39/61/210/232/254) as "Recorded"
#     with their real attendance numbers
#   - fix a students-recorded count on row 83

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "Recorded By" column (G): flip "System, dnasr281@gmail.com" ->
#    "dnasr281@gmail.com, System" for every session row that lists both.
# ---------------------------------------------------------------------------
$gRows = @(2,3,4,23,24,25,26,45,46,47,48,67,68,69,70,89,90,91,110,111,112,131,132,133,152,153,154,173,174,175,194,195,196,197,216,217,218,219,238,239,240,241)
foreach ($r in $gRows) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# 2. Narrow column I ("Students") from width 14 to width 10.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 9.17

# ---------------------------------------------------------------------------
# 3. Top summary box (K/L columns): Recorded Sessions / Missing Sessions
#    counts, plus Coverage % / Average Attendance % text.
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 0
$ws.Range("L9").Value = "'30.2%"
$ws.Range("L10").Value = "'81.6%"

# ---------------------------------------------------------------------------
# 4. Per-group breakdown table (K:S columns) for the groups whose numbers
#    shifted: B1-10 (row 16), B1-11 (row 17), B1-7 (row 24), B1-8 (row 25),
#    B1-9 (row 26); plus B1-12 (row 18) whose attendance % alone moved.
# ---------------------------------------------------------------------------
$ws.Range("O16").Value = 7
$ws.Range("P16").Value = 0
$ws.Range("R16").Value = "'31.8%"
$ws.Range("S16").Value = "'78.8%"

$ws.Range("O17").Value = 7
$ws.Range("P17").Value = 0
$ws.Range("R17").Value = "'31.8%"
$ws.Range("S17").Value = "'66.9%"

$ws.Range("S18").Value = "'87.8%"

$ws.Range("O24").Value = 7
$ws.Range("P24").Value = 0
$ws.Range("R24").Value = "'31.8%"
$ws.Range("S24").Value = "'71.4%"

$ws.Range("O25").Value = 7
$ws.Range("P25").Value = 0
$ws.Range("R25").Value = "'31.8%"
$ws.Range("S25").Value = "'77.8%"

$ws.Range("O26").Value = 7
$ws.Range("P26").Value = 0
$ws.Range("R26").Value = "'31.8%"
$ws.Range("S26").Value = "'72.9%"

# ---------------------------------------------------------------------------
# 5. Row 83 (B1-12, session 2): attendance recount 2/21 -> 17/21.
# ---------------------------------------------------------------------------
$ws.Range("H83").Value = "17/21"

# ---------------------------------------------------------------------------
# 6. Sessions that were "Not Recorded" (pink highlight, style matching the
#    pending rows) are now "Recorded" (green highlight, same style as the
#    other recorded rows). Re-stamp the formatting from an already-recorded
#    row (16) and then fill in the real "Recorded By" / "Students" / "Status"
#    values for each of the five newly-recorded sessions.
# ---------------------------------------------------------------------------
$newlyRecorded = @{
    39  = @{ H = "23/31"; Students = 31 }
    61  = @{ H = "9/19";  Students = 19 }
    210 = @{ H = "16/27"; Students = 27 }
    232 = @{ H = "18/29"; Students = 29 }
    254 = @{ H = "18/29"; Students = 29 }
}

foreach ($r in $newlyRecorded.Keys) {
    $ws.Range("A16:I16").Copy()
    $ws.Range("A$r`:I$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 7).Value = "System"
    $ws.Cells.Item($r, 8).Value = $newlyRecorded[$r].H
    $ws.Cells.Item($r, 9).Value = "Recorded"
}

$excel.CutCopyMode = 0
